# Got to algo 5 with a little cleanup.
#
# 1) Duplicate "Algorithm 4" into a new trailing sheet "Algorithm 5"
#    (this carries over formatting/styles exactly, matching how the
#    original author appears to have built each new algorithm sheet).
# 2) Update the text on the new sheet to describe algorithm 5 (splitting
#    a scholarship amount across applicants) instead of algorithm 4's
#    "get top ranked applicant who has not already won" logic.
# 3) Move the active selection to H19 on the new sheet.

$wb = $excel.ActiveWorkbook

$sheet4 = $wb.Worksheets.Item("Algorithm 4")
$sheet4.Copy([System.Reflection.Missing]::Value, $sheet4)

$ws = $wb.Worksheets.Item($sheet4.Index + 1)
$ws.Name = "Algorithm 5"

# --- Update the body text for the new algorithm ------------------------
# (values are entered in the same order the author appears to have typed
# them in, so new shared-string entries land in the same sequence)

$ws.Range("B9").Value = "Remove previous results from result grid"
$ws.Range("B10").Value = "Set Scholarship Loop Order By Scholarship Amount Descending"

$ws.Range("B11").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("B13").ClearContents()

$ws.Range("C12").Value = "Declare @CurrentAmount"
$ws.Range("C13").Value = "Declare @CurrentSplitAmount"
$ws.Range("C14").Value = "declare @CountOfApplicants"
$ws.Range("C11").Value = "Get Scholarship And Applicant Data"
$ws.Range("C15").Value = "Set @CurrentAmount and @CountOfApplicants"
$ws.Range("C16").Value = "set @Currentsplitamount as @CurrentAmount/@CountOfApplicants"
$ws.Range("C17").Value = "Foreach Applicants"

$ws.Range("C18").ClearContents()
$ws.Range("D18").Value = "Set Result with @currentsplitamount"

# Row 12 picked up a slight custom height in the authored file.
$ws.Rows.Item(12).RowHeight = 14.25

# --- Selection on the new sheet -----------------------------------------

$ws.Range("H19").Select()
